# Update Active_Outages.xlsx - 6/18/2025, 5:03:50 PM
# Applies elapsed-duration / hub-site / battery-backup-status refresh edits
# across sheets R1, R2, R4, R5 and R6.

$wb = $excel.ActiveWorkbook

$sheetR1 = $wb.Worksheets.Item("R1")
$sheetR1.Range("G2").Value = "3930:17:39"
$sheetR1.Range("G3").Value = "69:50:17"
$sheetR1.Range("G4").Value = "92:50:17"

$sheetR2 = $wb.Worksheets.Item("R2")
$sheetR2.Range("G2").Value = "12111:41:29"
$sheetR2.Range("G3").Value = "3241:24:58"
$sheetR2.Range("G4").Value = "479:36:32"
$sheetR2.Range("D5").Value = "JED0155"
$sheetR2.Range("J5").Value = "Good+Vandalized"

$sheetR4 = $wb.Worksheets.Item("R4")
$sheetR4.Range("G2").Value = "2957:31:09"
$sheetR4.Range("G3").Value = "184:43:24"
$sheetR4.Range("G4").Value = "72:55:49"
$sheetR4.Range("G5").Value = "70:33:22"

$sheetR5 = $wb.Worksheets.Item("R5")
$sheetR5.Range("G2").Value = "431:30:08"

$sheetR6 = $wb.Worksheets.Item("R6")
$sheetR6.Range("G2").Value = "72:02:26"
